$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph entirely (it sits right
#    after the H1 "Play Aurora Beast Hunter Free Today" heading).
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Meta description*") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) At the very end of the document, turn the single DALLE-prompt
#    paragraph into two paragraphs:
#      - a new bold paragraph: "Play Aurora Beast Hunter Free Today"
#      - the existing (italic) paragraph, but with its text replaced by
#        the meta-description sentence instead of the DALLE prompt.
# ------------------------------------------------------------------
$lastPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*DALLE*") {
        $lastPara = $p
        break
    }
}

$newXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Aurora Beast Hunter Free Today</w:t></w:r></w:p>
<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Experience big wins and exciting features in Aurora Beast Hunter, a 5-reel, 40-payline slot game from Just for the Win and Microgaming. Play free now.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

[void]$lastPara.Range.InsertXML($newXml)

# The XML insertion above replaces the old last paragraph's content with
# two brand-new paragraphs, but it also leaves a stray empty paragraph
# behind it (Word always needs a trailing paragraph mark for the body).
# Merge that stray empty paragraph away by deleting the paragraph break
# directly before it.
$newCount = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($newCount - 1)
$trailing = $d.Paragraphs.Item($newCount)
$mergeRange = $d.Range($secondLast.Range.End - 1, $trailing.Range.End)
$mergeRange.Delete()
